$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7.5619187107448
$ws.Range("D2").Value = 14.94569500033405
$ws.Range("E2").Value = 8.133560871811834
$ws.Range("F2").Value = 108.2532446164392
$ws.Range("G2").Value = 4.134971777197348
$ws.Range("J2").Value = 11.18115815630198
$ws.Range("L2").Value = 9.606543582882365
$ws.Range("M2").Value = 75.39830013714037
$ws.Range("C3").Value = 7.626021159278602
$ws.Range("D3").Value = 15.01360064984741
$ws.Range("E3").Value = 7.79799673708996
$ws.Range("F3").Value = 109.5234281847744
$ws.Range("G3").Value = 4.151887848442145
$ws.Range("J3").Value = 11.31151842423965
$ws.Range("L3").Value = 9.483884824707939
$ws.Range("M3").Value = 73.5573528289397
$ws.Range("C4").Value = 7.66700799315175
$ws.Range("D4").Value = 15.06444122107265
$ws.Range("E4").Value = 7.586345404946105
$ws.Range("F4").Value = 110.3592733583721
$ws.Range("G4").Value = 4.162662434652582
$ws.Range("J4").Value = 11.39458764928221
$ws.Range("L4").Value = 9.409781009756063
$ws.Range("M4").Value = 72.41580856144762
$ws.Range("C5").Value = 7.684124256172367
$ws.Range("D5").Value = 15.08739660221239
$ws.Range("E5").Value = 7.498786847267606
$ws.Range("F5").Value = 110.7137077342582
$ws.Range("G5").Value = 4.167152507713094
$ws.Range("J5").Value = 11.42921004646317
$ws.Range("L5").Value = 9.37990413033976
$ws.Range("M5").Value = 71.94826033604933
$ws.Range("C6").Value = 7.686991520795708
$ws.Range("D6").Value = 15.09134165959729
$ws.Range("E6").Value = 7.484172033425753
$ws.Range("F6").Value = 110.7733891357339
$ws.Range("G6").Value = 4.167904131195297
$ws.Range("J6").Value = 11.43500591803291
$ws.Range("L6").Value = 9.374962979363758
$ws.Range("M6").Value = 71.87049562703554
$ws.Range("C7").Value = 7.667237147714859
$ws.Range("D7").Value = 15.0647418276352
$ws.Range("E7").Value = 7.585169713506335
$ws.Range("F7").Value = 110.3639977384074
$ws.Range("G7").Value = 4.16272258483132
$ws.Range("J7").Value = 11.39505144445413
$ws.Range("L7").Value = 9.40937675756668
$ws.Range("M7").Value = 72.40951197512466
$ws.Range("C8").Value = 7.583686159413825
$ws.Range("D8").Value = 14.96717260378404
$ws.Range("E8").Value = 8.019080765299893
$ws.Range("F8").Value = 108.6794372236339
$ws.Range("G8").Value = 4.140724983574929
$ws.Range("J8").Value = 11.22548393585497
$ws.Range("L8").Value = 9.564010072326099
$ws.Range("M8").Value = 74.76612946958124
$ws.Range("C9").Value = 7.432547918264095
$ws.Range("D9").Value = 14.85127285797558
$ws.Range("E9").Value = 8.821628100309066
$ws.Range("F9").Value = 105.8320543815007
$ws.Range("G9").Value = 4.100583075823963
$ws.Range("J9").Value = 10.91653451472366
$ws.Range("L9").Value = 9.876104644826023
$ws.Range("M9").Value = 79.28019278644621
$ws.Range("C10").Value = 7.328950842782504
$ws.Range("D10").Value = 14.81613021766198
$ws.Range("E10").Value = 9.377607491418258
$ws.Range("F10").Value = 104.0364737141251
$ws.Range("G10").Value = 4.072794115713639
$ws.Range("J10").Value = 10.70328024406949
$ws.Range("L10").Value = 10.1098536788017
$ws.Range("M10").Value = 82.50910830645513
$ws.Range("C11").Value = 7.283370404883856
$ws.Range("D11").Value = 14.81189519658585
$ws.Range("E11").Value = 9.622587655467292
$ws.Range("F11").Value = 103.2885316750018
$ws.Range("G11").Value = 4.06049334753298
$ws.Range("J11").Value = 10.60910514975649
$ws.Range("L11").Value = 10.21697159319924
$ws.Range("M11").Value = 83.95506440407586
$ws.Range("C12").Value = 7.266326971552097
$ws.Range("D12").Value = 14.81205681453685
$ws.Range("E12").Value = 9.714174181370536
$ws.Range("F12").Value = 103.0156348788872
$ws.Range("G12").Value = 4.055881834473954
$ws.Range("J12").Value = 10.57383957580459
$ws.Range("L12").Value = 10.25763244768089
$ws.Range("M12").Value = 84.49903906602935
$ws.Range("C13").Value = 7.269988023052348
$ws.Range("D13").Value = 14.81194227529496
$ws.Range("E13").Value = 9.694502507436214
$ws.Range("F13").Value = 103.0739413051266
$ws.Range("G13").Value = 4.05687297799203
$ws.Range("J13").Value = 10.58141719531355
$ws.Range("L13").Value = 10.24887129085202
$ws.Range("M13").Value = 84.38204761900018
$ws.Range("C14").Value = 7.281963916677875
$ws.Range("D14").Value = 14.81187267529959
$ws.Range("E14").Value = 9.630146425260198
$ws.Range("F14").Value = 103.2658706186837
$ws.Range("G14").Value = 4.060113037824658
$ws.Range("J14").Value = 10.60619594861577
$ws.Range("L14").Value = 10.22031491187907
$ws.Range("M14").Value = 83.99989071827659
$ws.Range("C15").Value = 7.289327576708754
$ws.Range("D15").Value = 14.81206220782789
$ws.Range("E15").Value = 9.590571471675839
$ws.Range("F15").Value = 103.3847918499413
$ws.Range("G15").Value = 4.062103651026758
$ws.Range("J15").Value = 10.62142494822069
$ws.Range("L15").Value = 10.20283559068902
$ws.Range("M15").Value = 83.7653348947251
$ws.Range("C16").Value = 7.331960284973387
$ws.Range("D16").Value = 14.8166498731547
$ws.Range("E16").Value = 9.36143403340183
$ws.Range("F16").Value = 104.0867751935512
$ws.Range("G16").Value = 4.073604646750114
$ws.Range("J16").Value = 10.70949090759069
$ws.Range("L16").Value = 10.10286772946827
$ws.Range("M16").Value = 82.41412568132259
$ws.Range("C17").Value = 7.3585063009035
$ws.Range("D17").Value = 14.82252731358638
$ws.Range("E17").Value = 9.218800536440959
$ws.Range("F17").Value = 104.5353528998292
$ws.Range("G17").Value = 4.080745665108861
$ws.Range("J17").Value = 10.76423480247323
$ws.Range("L17").Value = 10.04173039559663
$ws.Range("M17").Value = 81.57912320744427
$ws.Range("C18").Value = 7.37392076016582
$ws.Range("D18").Value = 14.82700936300668
$ws.Range("E18").Value = 9.136016181087335
$ws.Range("F18").Value = 104.7998180483481
$ws.Range("G18").Value = 4.084885201469073
$ws.Range("J18").Value = 10.79598960439661
$ws.Range("L18").Value = 10.00664047815697
$ws.Range("M18").Value = 81.09671156169867
$ws.Range("C19").Value = 7.379165052266313
$ws.Range("D19").Value = 14.82871381765943
$ws.Range("E19").Value = 9.107860208118259
$ws.Range("F19").Value = 104.8904576488005
$ws.Range("G19").Value = 4.086292381990359
$ws.Range("J19").Value = 10.80678752710216
$ws.Range("L19").Value = 9.994772939814219
$ws.Range("M19").Value = 80.93301701858472
$ws.Range("C20").Value = 7.355665371891719
$ws.Range("D20").Value = 14.82178710713519
$ws.Range("E20").Value = 9.234061581215837
$ws.Range("F20").Value = 104.486929977267
$ws.Range("G20").Value = 4.079982175035248
$ws.Range("J20").Value = 10.75837961043778
$ws.Range("L20").Value = 10.04823095589633
$ws.Range("M20").Value = 81.66823435465572
$ws.Range("C21").Value = 7.278440464768356
$ws.Range("D21").Value = 14.811844598745
$ws.Range("E21").Value = 9.649081726778439
$ws.Range("F21").Value = 103.2092122963528
$ws.Range("G21").Value = 4.059160111455741
$ws.Range("J21").Value = 10.59890715808797
$ws.Range("L21").Value = 10.2287000775256
$ws.Range("M21").Value = 84.11223863265882
$ws.Range("C22").Value = 7.229231032458018
$ws.Range("D22").Value = 14.81567242078626
$ws.Range("E22").Value = 9.913418646224699
$ws.Range("F22").Value = 102.4346034941342
$ws.Range("G22").Value = 4.045821597748471
$ws.Range("J22").Value = 10.49698837300797
$ws.Range("L22").Value = 10.34721004090475
$ws.Range("M22").Value = 85.68858019383975
$ws.Range("C23").Value = 7.255381467503795
$ws.Range("D23").Value = 14.81265928218134
$ws.Range("E23").Value = 9.772979671586139
$ws.Range("F23").Value = 102.842343517883
$ws.Range("G23").Value = 4.052916767892761
$ws.Range("J23").Value = 10.5511771443593
$ws.Range("L23").Value = 10.28391219684453
$ws.Range("M23").Value = 84.84925845551179
$ws.Range("C24").Value = 7.356949279544741
$ws.Range("D24").Value = 14.822118324473
$ws.Range("E24").Value = 9.227164496858839
$ws.Range("F24").Value = 104.5088015426342
$ws.Range("G24").Value = 4.080327242594615
$ws.Range("J24").Value = 10.76102586414224
$ws.Range("L24").Value = 10.04529186842435
$ws.Range("M24").Value = 81.62795452627306
$ws.Range("C25").Value = 7.472105226142542
$ws.Range("D25").Value = 14.87414131481847
$ws.Range("E25").Value = 8.610172406607823
$ws.Range("F25").Value = 106.5517309988615
$ws.Range("G25").Value = 4.111134139016536
$ws.Range("J25").Value = 10.99765523689497
$ws.Range("L25").Value = 9.790826006683293
$ws.Range("M25").Value = 78.07292796256394
